$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Totale" label (D19) to "Totale fallimenti"
$ws.Range("D19").Value = "Totale fallimenti"

# Row 19 header labels for the three new columns
$ws.Range("E19").Value = "Generazioni mancanti"
$ws.Range("F19").Value = "Generazioni non necessarie"
$ws.Range("G19").Value = "Generazioni necessarie ma errate"

# Row 20 values (LLM totals)
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 0

# Row 21 values (Analitica totals)
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0

# Column width adjustments (columns F and G got wider to fit the new headers)
$ws.Columns.Item(6).ColumnWidth = 26.333333333333332
$ws.Columns.Item(7).ColumnWidth = 29

# Update the view: move selection to G21 (this also drops the stale topLeftCell scroll anchor)
$ws.Range("G21").Select()
